# Apply the "O" marks for week 1 (column F, row 5 -> 네트워크/5주차? actually 5th course row)
# and row 8 to the timetable, matching the xl/worksheets/sheet1.xml diff:
#   F5: "" -> "O" (shared string idx 15)
#   F8: "" -> "O" (shared string idx 15)
# and move the active selection from E10 to F9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = "O"
$ws.Range("F8").Value = "O"

# Update the sheet's saved selection/active cell (sheetView/selection element).
$ws.Range("F9").Select()
